# Aggiornamento fino a 13/03 - append 4 new daily rows (252-255) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (style "s=2": center/top aligned, thin
# border, custom date/time number format) from the last existing row
# down onto the 4 new date cells before filling in values, so the new
# rows render identically to the rest of column A.
$ws.Range("A251").Copy() | Out-Null
$ws.Range("A252:A255").PasteSpecial(-4122) | Out-Null

$ws.Range("A252").Value = 44326
$ws.Range("B252").Value = 1
$ws.Range("C252").Value = 26
$ws.Range("D252").Value = 147.9879332915932

$ws.Range("A253").Value = 44327
$ws.Range("B253").Value = 5
$ws.Range("C253").Value = 30
$ws.Range("D253").Value = 170.7553076441459

$ws.Range("A254").Value = 44328
$ws.Range("B254").Value = 0
$ws.Range("C254").Value = 30
$ws.Range("D254").Value = 170.7553076441459

$ws.Range("A255").Value = 44329
$ws.Range("B255").Value = 1
$ws.Range("C255").Value = 23
$ws.Range("D255").Value = 130.9124025271786
